# Fixing a host of issues with select_one_with_other
#
# Adds a new survey question (select_one_with_other colors) plus a note
# referencing a selected() call whose argument isn't in the choices list,
# and adds the backing "colors" choice list (red/green/blue).

$wb = $excel.ActiveWorkbook

# --- survey sheet: add the new question + note rows ---
$survey = $wb.Worksheets.Item("survey")

$survey.Range("B14").Value = "select_one_with_other colors"
$survey.Range("E14").Value = "color"
$survey.Range("F14").Value = "What is your favorite color?"

$survey.Range("A15").Value = "selected function with arguement not included in choices."
$survey.Range("B15").Value = "note"
$survey.Range("D15").Value = 'selected(data(''color''), ''teal'')'
$survey.Range("F15").Value = "Teal is a good choice."

# --- choices sheet: add the new "colors" choice list ---
$choices = $wb.Worksheets.Item("choices")

$choices.Range("A15").Value = "colors"
$choices.Range("B15").Value = "red"
$choices.Range("D15").Value = "Red"

$choices.Range("A16").Value = "colors"
$choices.Range("B16").Value = "green"
$choices.Range("D16").Value = "Green"

$choices.Range("A17").Value = "colors"
$choices.Range("B17").Value = "blue"
$choices.Range("D17").Value = "Blue"
